$d = $word.ActiveDocument

# Hunk 1: paragraph 2 (paraId 4301507D) - empty paragraph becomes the bold
# "4조 김상규, 손민우, 최희영" line with 8 leading tabs, plus the _GoBack bookmark
# (the bookmark used to live at the end of paragraph 3; Word moved it here).
$p1 = $d.Paragraphs.Item(2)
$r1 = $p1.Range
$xml1 = '<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4301507D" w14:textId="77777777" w:rsidR="00A553DA" w:rsidRDefault="00A553DA" w:rsidP="00A553DA"><w:pPr><w:pStyle w:val="a3"/><w:ind w:leftChars="0" w:left="760"/><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/></w:rPr></w:pPr><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>4</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/></w:rPr><w:t>조 김상규,</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/></w:rPr><w:t>손민우,</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/></w:rPr><w:t>최희영</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r1.InsertXML($xml1)

# Hunk 2: paragraph 3 (paraId 2479954B) - drop the _GoBack bookmark that used
# to sit at the end of this paragraph (now on paragraph 2 instead).
$p2 = $d.Paragraphs.Item(3)
$r2 = $p2.Range
$xml2 = '<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2479954B" w14:textId="77777777" w:rsidR="00D13964" w:rsidRPr="00A553DA" w:rsidRDefault="00D13964" w:rsidP="00D13964"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00A553DA"><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Standard, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00A553DA"><w:rPr><w:sz w:val="24"/></w:rPr><w:t>MinMax</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00A553DA"><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">, Max-Abs, Robust </w:t></w:r><w:r w:rsidRPr="00A553DA"><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="24"/></w:rPr><w:t>함수 구현 및 설명</w:t></w:r></w:p><w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r2.InsertXML($xml2)

# Hunk 3: paragraph 24 (paraId 341B9011, last paragraph before sectPr) - the
# paragraph-mark run properties no longer carry the eastAsia font hint.
$p3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$r3 = $p3.Range
$xml3 = '<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="341B9011" w14:textId="77777777" w:rsidR="00D13964" w:rsidRDefault="00D13964" w:rsidP="00D13964"><w:pPr><w:ind w:left="400"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="2E1CDA47" wp14:editId="23EEA9DD"><wp:extent cx="5731510" cy="3350895"/><wp:effectExtent l="0" t="0" r="2540" b="0"/><wp:docPr id="8" name="그림 8" descr="C:\Users\wowze\AppData\Local\Microsoft\Windows\INetCache\Content.MSO\9C640992.tmp"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 7" descr="C:\Users\wowze\AppData\Local\Microsoft\Windows\INetCache\Content.MSO\9C640992.tmp"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId12"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="5731510" cy="3350895"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p><w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r3.InsertXML($xml3)

Write-Host "all hunks applied"
